$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation for Price column (D) updates so values
# like "315.59" are not silently coerced into numbers by Excel.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '28.404.63'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.819.93'
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '315.59'
$ws.Range('E5').Value = '  -0.62%  '
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5134'
$ws.Range('E7').Value = '  -3.71%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3918'
$ws.Range('E8').Value = '  -4.02%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07864'
$ws.Range('E9').Value = '  +3.50%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '41.74'
$ws.Range('E10').Value = '  -0.24%  '
$ws.Range('E11').Value = '  +0.58%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '20.94'
$ws.Range('E12').Value = '  +0.96%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.241'
$ws.Range('E13').Value = '  -1.33%  '
$ws.Range('E14').Value = '  -0.18%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.476'
$ws.Range('E15').Value = '  -1.09%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.822.18'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001127'
$ws.Range('E17').Value = '  +5.08%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '92.50'
$ws.Range('E18').Value = '  +3.68%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06623'
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.67'
$ws.Range('E20').Value = '  +0.54%  '
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.080'
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '28.433.26'
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('E24').Value = '  +0.82%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.266'
$ws.Range('E25').Value = '  +3.88%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '21.08'
$ws.Range('E26').Value = '  +2.64%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.028.06'
$ws.Range('E27').Value = '  -0.61%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '154.32'
$ws.Range('E28').Value = '  -1.46%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.395'
$ws.Range('E29').Value = '  -2.57%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '125.40'
$ws.Range('E30').Value = '  +1.48%  '
$ws.Range('E31').Value = '  +0.70%  '
$ws.Range('E32').Value = '  -1.29%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.668'
$ws.Range('E33').Value = '  +0.17%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.647'
$ws.Range('E34').Value = '  +0.46%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.07037'
$ws.Range('E35').Value = '  -2.00%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.2211'
$ws.Range('E36').Value = '  -1.91%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02322'
$ws.Range('E37').Value = '  -0.49%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '5.193'
$ws.Range('E38').Value = '  -0.24%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '8.784'
$ws.Range('E39').Value = '  -0.42%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.6250'
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '11.23'
$ws.Range('E41').Value = '  -0.42%  '
$ws.Range('E42').Value = '  -0.77%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.000'
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('E44').Value = '  -0.83%  '
$ws.Range('E45').Value = '  -0.23%  '
$ws.Range('E46').Value = '  +0.78%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5881'
$ws.Range('E47').Value = '  +0.81%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '124.57'
$ws.Range('E48').Value = '  -0.80%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.970'
$ws.Range('E49').Value = '  -0.60%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.192'
$ws.Range('E50').Value = '  -1.49%  '
$ws.Range('E51').Value = '  -0.03%  '
